$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column before column B to make room for "Period"
$ws.Columns("B:B").Insert()

# Header
$ws.Range("B1").Value = "Period"

# Data values
$ws.Range("B2").Value = "Quarter"
$ws.Range("B3").Value = "Quarter"
$ws.Range("B4").Value = "Quarter"
$ws.Range("B5").Value = "Annual"

# Match the column width of the new Period column to column A
$ws.Columns("B:B").ColumnWidth = $ws.Columns("A:A").ColumnWidth

# Update the active selection to match the target state
$ws.Range("B6").Select()
